# Auto-generated edit script: append 12 new survey response rows (434-445)
# to the "Form_Responses1" table, matching rows scraped from a Google Form export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Template rows already present in the table, used to replicate the banded-row
# border/number-format styling (even rows, odd rows, and the special last-row style)
# that Excel applies automatically to table rows.
$evenTemplate = $ws.Range("A432:Y432")   # style used by even data rows (s=4/5/6)
$oddTemplate  = $ws.Range("A431:Y431")   # style used by odd data rows (s=7/8/9)
$lastTemplate = $ws.Range("A433:Y433")   # style used by the tables last row (s=16/17/18)

# ---- row 434 ----
$evenTemplate.Copy($ws.Range("A434:Y434"))
$ws.Rows.Item(434).RowHeight = 15.75
$ws.Range("A434").Value = 45610.880750902783
$ws.Range("B434").Value = "yeel6945@naver.com"
$ws.Range("C434").Value = "경영학과"
$ws.Range("D434").Value = 20203635
$ws.Range("E434").Value = "이수빈"
$ws.Range("F434").Value = "3. 3-5일"
$ws.Range("G434").Value = "3. 3-5일"
$ws.Range("H434").Value = "3. 3-5일"
$ws.Range("I434").Value = "3. 3-5일"
$ws.Range("J434").Value = "5. 6-7일"
$ws.Range("K434").Value = "5. 6-7일"
$ws.Range("L434").Value = "5. 6-7일"
$ws.Range("M434").Value = "5. 6-7일"
$ws.Range("N434").Value = "5. 6-7일"
$ws.Range("O434").Value = "5. 6-7일"
$ws.Range("P434").Value = "3. 가끔"
$ws.Range("Q434").Value = "3. 가끔"
$ws.Range("R434").Value = "3. 가끔"
$ws.Range("S434").Value = "3. 가끔"
$ws.Range("T434").Value = "5. 아니오"
$ws.Range("U434").Value = "5. 아니오"
$ws.Range("V434").Value = "5. 아니오"
$ws.Range("W434").Value = "3. 가끔"
$ws.Range("X434").Value = "3. 가끔"
$ws.Range("Y434").Value = "3. 가끔"

# ---- row 435 ----
$oddTemplate.Copy($ws.Range("A435:Y435"))
$ws.Rows.Item(435).RowHeight = 15.75
$ws.Range("A435").Value = 45610.961117245373
$ws.Range("B435").Value = "tngusvhs@gmail.com"
$ws.Range("C435").Value = "생명과학과"
$ws.Range("D435").Value = 20243529
$ws.Range("E435").Value = "이수현"
$ws.Range("F435").Value = "3. 3-5일"
$ws.Range("G435").Value = "5. 6-7일"
$ws.Range("H435").Value = "5. 6-7일"
$ws.Range("I435").Value = "3. 3-5일"
$ws.Range("J435").Value = "5. 6-7일"
$ws.Range("K435").Value = "5. 6-7일"
$ws.Range("L435").Value = "5. 6-7일"
$ws.Range("M435").Value = "5. 6-7일"
$ws.Range("N435").Value = "5. 6-7일"
$ws.Range("O435").Value = "3. 3-5일"
$ws.Range("P435").Value = "3. 가끔"
$ws.Range("Q435").Value = "3. 가끔"
$ws.Range("R435").Value = "5. 아니오"
$ws.Range("S435").Value = "5. 아니오"
$ws.Range("T435").Value = "3. 가끔"
$ws.Range("U435").Value = "3. 가끔"
$ws.Range("V435").Value = "3. 가끔"
$ws.Range("W435").Value = "5. 아니오"
$ws.Range("X435").Value = "5. 아니오"
$ws.Range("Y435").Value = "3. 가끔"

# ---- row 436 ----
$evenTemplate.Copy($ws.Range("A436:Y436"))
$ws.Rows.Item(436).RowHeight = 15.75
$ws.Range("A436").Value = 45610.963458090278
$ws.Range("B436").Value = "jerryterryharry@gmail.com"
$ws.Range("C436").Value = "빅데이터"
$ws.Range("D436").Value = 20205162
$ws.Range("E436").Value = "문진영"
$ws.Range("F436").Value = "1. 0-2일"
$ws.Range("G436").Value = "5. 6-7일"
$ws.Range("H436").Value = "3. 3-5일"
$ws.Range("I436").Value = "3. 3-5일"
$ws.Range("J436").Value = "3. 3-5일"
$ws.Range("K436").Value = "1. 0-2일"
$ws.Range("L436").Value = "1. 0-2일"
$ws.Range("M436").Value = "5. 6-7일"
$ws.Range("N436").Value = "5. 6-7일"
$ws.Range("O436").Value = "1. 0-2일"
$ws.Range("P436").Value = "3. 가끔"
$ws.Range("Q436").Value = "1. 예"
$ws.Range("R436").Value = "3. 가끔"
$ws.Range("S436").Value = "3. 가끔"
$ws.Range("T436").Value = "5. 아니오"
$ws.Range("U436").Value = "5. 아니오"
$ws.Range("V436").Value = "3. 가끔"
$ws.Range("W436").Value = "3. 가끔"
$ws.Range("X436").Value = "5. 아니오"
$ws.Range("Y436").Value = "5. 아니오"

# ---- row 437 ----
$oddTemplate.Copy($ws.Range("A437:Y437"))
$ws.Rows.Item(437).RowHeight = 15.75
$ws.Range("A437").Value = 45610.996303356485
$ws.Range("B437").Value = "jb9517asd@naver.com"
$ws.Range("C437").Value = "소프트웨어학부"
$ws.Range("D437").Value = 20245109
$ws.Range("E437").Value = "곽우주"
$ws.Range("F437").Value = "5. 6-7일"
$ws.Range("G437").Value = "5. 6-7일"
$ws.Range("H437").Value = "5. 6-7일"
$ws.Range("I437").Value = "3. 3-5일"
$ws.Range("J437").Value = "3. 3-5일"
$ws.Range("K437").Value = "3. 3-5일"
$ws.Range("L437").Value = "3. 3-5일"
$ws.Range("M437").Value = "5. 6-7일"
$ws.Range("N437").Value = "5. 6-7일"
$ws.Range("O437").Value = "5. 6-7일"
$ws.Range("P437").Value = "5. 아니오"
$ws.Range("Q437").Value = "5. 아니오"
$ws.Range("R437").Value = "3. 가끔"
$ws.Range("S437").Value = "5. 아니오"
$ws.Range("T437").Value = "5. 아니오"
$ws.Range("U437").Value = "5. 아니오"
$ws.Range("V437").Value = "5. 아니오"
$ws.Range("W437").Value = "5. 아니오"
$ws.Range("X437").Value = "5. 아니오"
$ws.Range("Y437").Value = "5. 아니오"

# ---- row 438 ----
$evenTemplate.Copy($ws.Range("A438:Y438"))
$ws.Rows.Item(438).RowHeight = 15.75
$ws.Range("A438").Value = 45611.038234560183
$ws.Range("B438").Value = "lhw2565@gmail.com"
$ws.Range("C438").Value = "미디어스쿨"
$ws.Range("D438").Value = 20242565
$ws.Range("E438").Value = "이혜원"
$ws.Range("F438").Value = "1. 0-2일"
$ws.Range("G438").Value = "5. 6-7일"
$ws.Range("H438").Value = "3. 3-5일"
$ws.Range("I438").Value = "3. 3-5일"
$ws.Range("J438").Value = "3. 3-5일"
$ws.Range("K438").Value = "5. 6-7일"
$ws.Range("L438").Value = "3. 3-5일"
$ws.Range("M438").Value = "5. 6-7일"
$ws.Range("N438").Value = "5. 6-7일"
$ws.Range("O438").Value = "1. 0-2일"
$ws.Range("P438").Value = "3. 가끔"
$ws.Range("Q438").Value = "3. 가끔"
$ws.Range("R438").Value = "5. 아니오"
$ws.Range("S438").Value = "5. 아니오"
$ws.Range("T438").Value = "5. 아니오"
$ws.Range("U438").Value = "3. 가끔"
$ws.Range("V438").Value = "3. 가끔"
$ws.Range("W438").Value = "3. 가끔"
$ws.Range("X438").Value = "5. 아니오"
$ws.Range("Y438").Value = "1. 예"

# ---- row 439 ----
$oddTemplate.Copy($ws.Range("A439:Y439"))
$ws.Rows.Item(439).RowHeight = 15.75
$ws.Range("A439").Value = 45611.390172141204
$ws.Range("B439").Value = "bigeyejimmy1@naver.com"
$ws.Range("C439").Value = "경영학과"
$ws.Range("D439").Value = 20182850
$ws.Range("E439").Value = "김현준"
$ws.Range("F439").Value = "1. 0-2일"
$ws.Range("G439").Value = "3. 3-5일"
$ws.Range("H439").Value = "3. 3-5일"
$ws.Range("I439").Value = "1. 0-2일"
$ws.Range("J439").Value = "1. 0-2일"
$ws.Range("K439").Value = "1. 0-2일"
$ws.Range("L439").Value = "1. 0-2일"
$ws.Range("M439").Value = "5. 6-7일"
$ws.Range("N439").Value = "5. 6-7일"
$ws.Range("O439").Value = "1. 0-2일"
$ws.Range("P439").Value = "5. 아니오"
$ws.Range("Q439").Value = "3. 가끔"
$ws.Range("R439").Value = "3. 가끔"
$ws.Range("S439").Value = "3. 가끔"
$ws.Range("T439").Value = "3. 가끔"
$ws.Range("U439").Value = "3. 가끔"
$ws.Range("V439").Value = "3. 가끔"
$ws.Range("W439").Value = "3. 가끔"
$ws.Range("X439").Value = "3. 가끔"
$ws.Range("Y439").Value = "5. 아니오"

# ---- row 440 ----
$evenTemplate.Copy($ws.Range("A440:Y440"))
$ws.Rows.Item(440).RowHeight = 15.75
$ws.Range("A440").Value = 45611.464371238428
$ws.Range("B440").Value = "yhh323@naver.com"
$ws.Range("C440").Value = "체육학과"
$ws.Range("D440").Value = 20184132
$ws.Range("E440").Value = "유형호"
$ws.Range("F440").Value = "3. 3-5일"
$ws.Range("G440").Value = "3. 3-5일"
$ws.Range("H440").Value = "3. 3-5일"
$ws.Range("I440").Value = "3. 3-5일"
$ws.Range("J440").Value = "3. 3-5일"
$ws.Range("K440").Value = "3. 3-5일"
$ws.Range("L440").Value = "1. 0-2일"
$ws.Range("M440").Value = "5. 6-7일"
$ws.Range("N440").Value = "5. 6-7일"
$ws.Range("O440").Value = "3. 3-5일"
$ws.Range("P440").Value = "3. 가끔"
$ws.Range("Q440").Value = "3. 가끔"
$ws.Range("R440").Value = "3. 가끔"
$ws.Range("S440").Value = "3. 가끔"
$ws.Range("T440").Value = "3. 가끔"
$ws.Range("U440").Value = "5. 아니오"
$ws.Range("V440").Value = "5. 아니오"
$ws.Range("W440").Value = "5. 아니오"
$ws.Range("X440").Value = "3. 가끔"
$ws.Range("Y440").Value = "5. 아니오"

# ---- row 441 ----
$oddTemplate.Copy($ws.Range("A441:Y441"))
$ws.Rows.Item(441).RowHeight = 15.75
$ws.Range("A441").Value = 45611.654391597222
$ws.Range("B441").Value = "hyj13223@naver.com"
$ws.Range("C441").Value = "정치행정학과"
$ws.Range("D441").Value = 20212432
$ws.Range("E441").Value = "이현진"
$ws.Range("F441").Value = "5. 6-7일"
$ws.Range("G441").Value = "5. 6-7일"
$ws.Range("H441").Value = "5. 6-7일"
$ws.Range("I441").Value = "5. 6-7일"
$ws.Range("J441").Value = "3. 3-5일"
$ws.Range("K441").Value = "3. 3-5일"
$ws.Range("L441").Value = "5. 6-7일"
$ws.Range("M441").Value = "5. 6-7일"
$ws.Range("N441").Value = "1. 0-2일"
$ws.Range("O441").Value = "5. 6-7일"
$ws.Range("P441").Value = "5. 아니오"
$ws.Range("Q441").Value = "5. 아니오"
$ws.Range("R441").Value = "5. 아니오"
$ws.Range("S441").Value = "5. 아니오"
$ws.Range("T441").Value = "5. 아니오"
$ws.Range("U441").Value = "5. 아니오"
$ws.Range("V441").Value = "5. 아니오"
$ws.Range("W441").Value = "5. 아니오"
$ws.Range("X441").Value = "5. 아니오"
$ws.Range("Y441").Value = "5. 아니오"

# ---- row 442 ----
$evenTemplate.Copy($ws.Range("A442:Y442"))
$ws.Rows.Item(442).RowHeight = 15.75
$ws.Range("A442").Value = 45611.696060891205
$ws.Range("B442").Value = "chaecjb@naver.com"
$ws.Range("C442").Value = "디지털미디어콘텐츠전공"
$ws.Range("D442").Value = 20203046
$ws.Range("E442").Value = "채희수"
$ws.Range("F442").Value = "5. 6-7일"
$ws.Range("G442").Value = "5. 6-7일"
$ws.Range("H442").Value = "5. 6-7일"
$ws.Range("I442").Value = "5. 6-7일"
$ws.Range("J442").Value = "5. 6-7일"
$ws.Range("K442").Value = "3. 3-5일"
$ws.Range("L442").Value = "5. 6-7일"
$ws.Range("M442").Value = "5. 6-7일"
$ws.Range("N442").Value = "5. 6-7일"
$ws.Range("O442").Value = "1. 0-2일"
$ws.Range("P442").Value = "3. 가끔"
$ws.Range("Q442").Value = "5. 아니오"
$ws.Range("R442").Value = "5. 아니오"
$ws.Range("S442").Value = "5. 아니오"
$ws.Range("T442").Value = "5. 아니오"
$ws.Range("U442").Value = "5. 아니오"
$ws.Range("V442").Value = "3. 가끔"
$ws.Range("W442").Value = "5. 아니오"
$ws.Range("X442").Value = "5. 아니오"
$ws.Range("Y442").Value = "1. 예"

# ---- row 443 ----
$oddTemplate.Copy($ws.Range("A443:Y443"))
$ws.Rows.Item(443).RowHeight = 15.75
$ws.Range("A443").Value = 45611.69865366898
$ws.Range("B443").Value = "bcy1976@naver.com"
$ws.Range("C443").Value = "빅데이터학과"
$ws.Range("D443").Value = 20235180
$ws.Range("E443").Value = "변치윤"
$ws.Range("F443").Value = "3. 3-5일"
$ws.Range("G443").Value = "3. 3-5일"
$ws.Range("H443").Value = "3. 3-5일"
$ws.Range("I443").Value = "3. 3-5일"
$ws.Range("J443").Value = "3. 3-5일"
$ws.Range("K443").Value = "3. 3-5일"
$ws.Range("L443").Value = "3. 3-5일"
$ws.Range("M443").Value = "3. 3-5일"
$ws.Range("N443").Value = "3. 3-5일"
$ws.Range("O443").Value = "3. 3-5일"
$ws.Range("P443").Value = "3. 가끔"
$ws.Range("Q443").Value = "3. 가끔"
$ws.Range("R443").Value = "3. 가끔"
$ws.Range("S443").Value = "3. 가끔"
$ws.Range("T443").Value = "3. 가끔"
$ws.Range("U443").Value = "3. 가끔"
$ws.Range("V443").Value = "3. 가끔"
$ws.Range("W443").Value = "3. 가끔"
$ws.Range("X443").Value = "3. 가끔"
$ws.Range("Y443").Value = "3. 가끔"

# ---- row 444 ----
$evenTemplate.Copy($ws.Range("A444:Y444"))
$ws.Rows.Item(444).RowHeight = 15.75
$ws.Range("A444").Value = 45611.708388067127
$ws.Range("B444").Value = "emf1811@naver.com"
$ws.Range("C444").Value = "바이오메디컬학과"
$ws.Range("D444").Value = 20233605
$ws.Range("E444").Value = "김들"
$ws.Range("F444").Value = "3. 3-5일"
$ws.Range("G444").Value = "5. 6-7일"
$ws.Range("H444").Value = "5. 6-7일"
$ws.Range("I444").Value = "5. 6-7일"
$ws.Range("J444").Value = "5. 6-7일"
$ws.Range("K444").Value = "5. 6-7일"
$ws.Range("L444").Value = "3. 3-5일"
$ws.Range("M444").Value = "5. 6-7일"
$ws.Range("N444").Value = "5. 6-7일"
$ws.Range("O444").Value = "5. 6-7일"
$ws.Range("P444").Value = "5. 아니오"
$ws.Range("Q444").Value = "1. 예"
$ws.Range("R444").Value = "5. 아니오"
$ws.Range("S444").Value = "5. 아니오"
$ws.Range("T444").Value = "5. 아니오"
$ws.Range("U444").Value = "5. 아니오"
$ws.Range("V444").Value = "5. 아니오"
$ws.Range("W444").Value = "5. 아니오"
$ws.Range("X444").Value = "3. 가끔"
$ws.Range("Y444").Value = "5. 아니오"

# ---- row 445 ----
$lastTemplate.Copy($ws.Range("A445:Y445"))
$ws.Rows.Item(445).RowHeight = 15.75
$ws.Range("A445").Value = 45611.723308703702
$ws.Range("B445").Value = "leyy2k@naver.com"
$ws.Range("C445").Value = "빅데이터"
$ws.Range("D445").Value = 202155115
$ws.Range("E445").Value = "김대현"
$ws.Range("F445").Value = "3. 3-5일"
$ws.Range("G445").Value = "3. 3-5일"
$ws.Range("H445").Value = "3. 3-5일"
$ws.Range("I445").Value = "3. 3-5일"
$ws.Range("J445").Value = "3. 3-5일"
$ws.Range("K445").Value = "3. 3-5일"
$ws.Range("L445").Value = "3. 3-5일"
$ws.Range("M445").Value = "3. 3-5일"
$ws.Range("N445").Value = "3. 3-5일"
$ws.Range("O445").Value = "3. 3-5일"
$ws.Range("P445").Value = "5. 아니오"
$ws.Range("Q445").Value = "3. 가끔"
$ws.Range("R445").Value = "5. 아니오"
$ws.Range("S445").Value = "5. 아니오"
$ws.Range("T445").Value = "5. 아니오"
$ws.Range("U445").Value = "5. 아니오"
$ws.Range("V445").Value = "5. 아니오"
$ws.Range("W445").Value = "3. 가끔"
$ws.Range("X445").Value = "5. 아니오"
$ws.Range("Y445").Value = "3. 가끔"

# Resize the table/ListObject to include the newly added rows
$tbl.Resize($ws.Range("A1:Y445"))

# Restore/update the view state (scroll position + active cell selection)
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 402
$ws.Range("C451").Select()

Write-Output "Added rows 434-445 to Form_Responses1 table"